$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("2:4").Delete()
$ws.Range("A11").Select()
